$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 17 (hunk 0)
$ws.Range("H17").Value = 1166.4412
$ws.Range("J17").Value = 1176.8507
$ws.Range("L17").Value = 3530.5521
$ws.Range("N17").Value = -3866.5521
# row 33 (hunk 1)
$ws.Range("H33").Value = 1487.0667
$ws.Range("I33").Value = 214
$ws.Range("K33").Value = 214
$ws.Range("M33").Value = 15
# row 86 (hunk 2)
$ws.Range("H86").Value = 2232.7727
$ws.Range("I86").Value = 2616.6667
$ws.Range("J86").Value = 1772.1
$ws.Range("K86").Value = 2616.6667
$ws.Range("L86").Value = 1772.1
$ws.Range("M86").Value = -1493.6667
$ws.Range("N86").Value = -4018.1
# row 88 (hunk 3)
$ws.Range("H88").Value = 5969.7
$ws.Range("J88").Value = 5969.7
$ws.Range("L88").Value = 5969.7
$ws.Range("N88").Value = -6781.7
# row 89 (hunk 4)
$ws.Range("H89").Value = 2232.7727
$ws.Range("I89").Value = 2616.6667
$ws.Range("J89").Value = 1772.1
$ws.Range("K89").Value = 13083.3335
$ws.Range("L89").Value = 8860.5
$ws.Range("M89").Value = -7467.333500000001
$ws.Range("N89").Value = -20092.5
# row 91 (hunk 5)
$ws.Range("H91").Value = 5969.7
$ws.Range("J91").Value = 5969.7
$ws.Range("L91").Value = 5969.7
$ws.Range("N91").Value = -8777.700000000001
# row 96 (hunk 6)
$ws.Range("H96").Value = 715087.4
$ws.Range("I96").Value = 833873.25
$ws.Range("J96").Value = 2372
$ws.Range("K96").Value = 2501619.75
$ws.Range("L96").Value = 7116
$ws.Range("M96").Value = -2500246.75
$ws.Range("N96").Value = -9862
# row 129 (hunk 7)
$ws.Range("H129").Value = 1607.8462
$ws.Range("I129").Value = 1162.3334
$ws.Range("J129").Value = 2610.25
$ws.Range("K129").Value = 3487.0002
$ws.Range("L129").Value = 7830.75
$ws.Range("M129").Value = 1512.9998
$ws.Range("N129").Value = -17830.75
# row 133 (hunk 8)
$ws.Range("H133").Value = 87500
$ws.Range("J133").Value = 87500
$ws.Range("L133").Value = 87500
$ws.Range("N133").Value = -97620
# row 134 (hunk 9)
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# row 135 (hunk 10)
$ws.Range("H135").Value = 1995.75
$ws.Range("I135").Value = 1997.8334
$ws.Range("J135").Value = 1989.5
$ws.Range("K135").Value = 17980.5006
$ws.Range("L135").Value = 17905.5
$ws.Range("M135").Value = -15445.5006
$ws.Range("N135").Value = -22975.5
# row 136 (hunk 11)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# row 137 (hunk 12)
$ws.Range("H137").Value = 1359.3077
$ws.Range("I137").Value = 1498.8667
$ws.Range("K137").Value = 4496.6001
$ws.Range("M137").Value = -1946.6001
# row 138 (hunk 13)
$ws.Range("H138").Value = 1710.7119
$ws.Range("I138").Value = 1568.25
$ws.Range("J138").Value = 1733.0588
$ws.Range("K138").Value = 4704.75
$ws.Range("L138").Value = 5199.1764
$ws.Range("M138").Value = 435.25
$ws.Range("N138").Value = -15479.1764
# row 141 (hunk 14)
$ws.Range("H141").Value = 3109.3076
$ws.Range("I141").Value = 2811.1365
$ws.Range("K141").Value = 8433.4095
$ws.Range("M141").Value = -3253.4095

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 15)
$ws.Range("H32").Value = 6657.278
$ws.Range("I32").Value = 4997.0312
$ws.Range("J32").Value = 19939.25
$ws.Range("K32").Value = 4997.0312
$ws.Range("L32").Value = 19939.25
$ws.Range("M32").Value = -4710.0312
$ws.Range("N32").Value = -20513.25
# row 39 (hunk 16)
$ws.Range("H39").Value = 675000
$ws.Range("I39").Value = 1002500
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 1002500
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = -1001980
$ws.Range("N39").Value = -21040
# row 51 (hunk 17)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# row 61 (hunk 18)
$ws.Range("H61").Value = 7535.9165
$ws.Range("I61").Value = 7492.6665
$ws.Range("J61").Value = 7665.6665
$ws.Range("K61").Value = 7492.6665
$ws.Range("L61").Value = 7665.6665
$ws.Range("M61").Value = -7280.6665
$ws.Range("N61").Value = -8089.6665
# row 74 (hunk 19)
$ws.Range("H74").Value = 2389.2354
$ws.Range("I74").Value = 2508
$ws.Range("J74").Value = 1498.5
$ws.Range("K74").Value = 2508
$ws.Range("L74").Value = 1498.5
$ws.Range("M74").Value = -1634
$ws.Range("N74").Value = -3246.5
# row 77 (hunk 20)
$ws.Range("H77").Value = 2389.2354
$ws.Range("I77").Value = 2508
$ws.Range("J77").Value = 1498.5
$ws.Range("K77").Value = 12540
$ws.Range("L77").Value = 7492.5
$ws.Range("M77").Value = -8172
$ws.Range("N77").Value = -16228.5
# row 110 (hunk 21)
$ws.Range("H110").Value = 253.77777
$ws.Range("I110").Value = 204.85715
$ws.Range("J110").Value = 425
$ws.Range("K110").Value = 204.85715
$ws.Range("L110").Value = 425
$ws.Range("M110").Value = 1840.14285
$ws.Range("N110").Value = -4515
# row 132 (hunk 22)
$ws.Range("H132").Value = 1821.0555
$ws.Range("I132").Value = 1594.9286
$ws.Range("K132").Value = 4784.7858
$ws.Range("M132").Value = -2254.7858
# row 136 (hunk 23)
$ws.Range("H136").Value = 7535.9165
$ws.Range("I136").Value = 7492.6665
$ws.Range("J136").Value = 7665.6665
$ws.Range("K136").Value = 22477.9995
$ws.Range("L136").Value = 22996.9995
$ws.Range("M136").Value = -19927.9995
$ws.Range("N136").Value = -28096.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 61 (hunk 24)
$ws.Range("H61").Value = 75000
$ws.Range("J61").Value = 75000
$ws.Range("L61").Value = 75000
$ws.Range("N61").Value = -75626
# row 62 (hunk 25)
$ws.Range("H62").Value = 43333.332
$ws.Range("J62").Value = 43333.332
$ws.Range("L62").Value = 43333.332
$ws.Range("N62").Value = -44705.332
# row 63 (hunk 26)
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
# row 65 (hunk 27)
$ws.Range("H65").Value = 43333.332
$ws.Range("J65").Value = 43333.332
$ws.Range("L65").Value = 129999.996
$ws.Range("N65").Value = -136863.996
# row 66 (hunk 28)
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
# row 134 (hunk 29)
$ws.Range("H134").Value = 2737.889
$ws.Range("I134").Value = 2790.3333
$ws.Range("K134").Value = 8370.999899999999
$ws.Range("M134").Value = -5835.999899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16 (hunk 30)
$ws.Range("H16").Value = 2328.7778
$ws.Range("I16").Value = 745
$ws.Range("K16").Value = 745
$ws.Range("M16").Value = -458
# row 31 (hunk 31)
$ws.Range("H31").Value = 3426.1875
$ws.Range("I31").Value = 1263.5
$ws.Range("K31").Value = 1263.5
$ws.Range("M31").Value = -968.5
# row 34 (hunk 32)
$ws.Range("H34").Value = 3426.1875
$ws.Range("I34").Value = 1263.5
$ws.Range("K34").Value = 1263.5
$ws.Range("M34").Value = -1061.5
# row 35 (hunk 33)
$ws.Range("H35").Value = 221.14285
$ws.Range("I35").Value = 174.66667
$ws.Range("K35").Value = 174.66667
$ws.Range("M35").Value = 119.33333
# row 113 (hunk 34)
$ws.Range("H113").Value = 2328.7778
$ws.Range("I113").Value = 745
$ws.Range("K113").Value = 745
$ws.Range("M113").Value = 1425
# row 132 (hunk 35)
$ws.Range("H132").Value = 1390.0264
$ws.Range("I132").Value = 1312.8529
$ws.Range("K132").Value = 3938.5587
$ws.Range("M132").Value = -1408.5587
# row 134 (hunk 36)
$ws.Range("H134").Value = 2422.675
$ws.Range("I134").Value = 1719.7222
$ws.Range("K134").Value = 5159.1666
$ws.Range("M134").Value = -2624.1666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 16 (hunk 37)
$ws.Range("H16").Value = 975
$ws.Range("I16").Value = 975
$ws.Range("K16").Value = 2925
$ws.Range("M16").Value = -2752
# row 33 (hunk 38)
$ws.Range("H33").Value = 163.375
$ws.Range("J33").Value = 215
$ws.Range("L33").Value = 1290
$ws.Range("N33").Value = -1856
# row 62 (hunk 39)
$ws.Range("H62").Value = 4850
$ws.Range("I62").Value = 4850
$ws.Range("K62").Value = 14550
$ws.Range("M62").Value = -13864
# row 65 (hunk 40)
$ws.Range("H65").Value = 4850
$ws.Range("I65").Value = 4850
$ws.Range("K65").Value = 43650
$ws.Range("M65").Value = -40218
# row 121 (hunk 41)
$ws.Range("H121").Value = 687.05554
$ws.Range("I121").Value = 278
$ws.Range("K121").Value = 834
$ws.Range("M121").Value = 476

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 119 (hunk 42)
$ws.Range("H119").Value = 72210
$ws.Range("J119").Value = 72210
$ws.Range("L119").Value = 72210
$ws.Range("N119").Value = -81886
# row 132 (hunk 43)
$ws.Range("H132").Value = 2805
$ws.Range("I132").Value = 2805
$ws.Range("K132").Value = 8415
$ws.Range("M132").Value = -5885

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40 (hunk 44)
$ws.Range("H40").Value = 5796.1943
$ws.Range("I40").Value = 3350.6843
$ws.Range("J40").Value = 8529.412
$ws.Range("K40").Value = 3350.6843
$ws.Range("L40").Value = 8529.412
$ws.Range("M40").Value = -3214.6843
$ws.Range("N40").Value = -8801.412
# row 100 (hunk 45)
$ws.Range("H100").Value = 4988.5186
$ws.Range("I100").Value = 3172.7334
$ws.Range("K100").Value = 3172.7334
$ws.Range("M100").Value = -2631.7334
# row 132 (hunk 46)
$ws.Range("H132").Value = 4112.773
$ws.Range("I132").Value = 3932.2778
$ws.Range("K132").Value = 11796.8334
$ws.Range("M132").Value = -9266.8334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 126 (hunk 47)
$ws.Range("H126").Value = 1849.0834
$ws.Range("I126").Value = 1521
$ws.Range("K126").Value = 4563
$ws.Range("M126").Value = -2093
# row 132 (hunk 48)
$ws.Range("H132").Value = 2203.6453
$ws.Range("I132").Value = 2416.3
$ws.Range("J132").Value = 1817
$ws.Range("K132").Value = 7248.900000000001
$ws.Range("L132").Value = 5451
$ws.Range("M132").Value = -4718.900000000001
$ws.Range("N132").Value = -10511
